$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.401.14"
$ws.Range("E2").Value = '  +1.58%  '
$ws.Range("D3").Value = "'1.940.43"
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = "'243.33"
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("E6").Value = '  -1.19%  '
$ws.Range("D7").Value = "'57.79"
$ws.Range("E7").Value = '  -5.47%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").Value = "'0.366"
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").Value = "'55.64"
$ws.Range("E10").Value = '  -0.52%  '
$ws.Range("D11").Value = "'0.0835"
$ws.Range("E11").Value = '  +5.80%  '
$ws.Range("E12").Value = '  +1.19%  '
$ws.Range("D13").Value = "'0.821"
$ws.Range("E13").Value = '  -3.48%  '
$ws.Range("D14").Value = "'21.45"
$ws.Range("E14").Value = '  -0.03%  '
$ws.Range("D15").Value = "'2.227.69"
$ws.Range("E15").Value = '  -0.20%  '
$ws.Range("E16").Value = '  -1.83%  '
$ws.Range("D17").Value = "'5.24"
$ws.Range("E17").Value = '  -2.75%  '
$ws.Range("D18").Value = "'1.941.05"
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").Value = "'36.373.47"
$ws.Range("E19").Value = '  +2.12%  '
$ws.Range("D20").Value = "'69.69"
$ws.Range("E20").Value = '  -0.98%  '
$ws.Range("D21").Value = "'0.0₃0863"
$ws.Range("E21").Value = '  +1.94%  '
$ws.Range("D22").Value = "'229.38"
$ws.Range("E22").Value = '  -3.66%  '
$ws.Range("D23").Value = "'5.04"
$ws.Range("E23").Value = '  -2.16%  '
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("E25").Value = '  -2.93%  '
$ws.Range("D26").Value = "'2.29"
$ws.Range("E26").Value = '  +0.64%  '
$ws.Range("D27").Value = "'9.24"
$ws.Range("E27").Value = '  -4.29%  '
$ws.Range("D28").Value = "'162.00"
$ws.Range("E28").Value = '  +2.52%  '
$ws.Range("D29").Value = "'19.37"
$ws.Range("E29").Value = '  -1.28%  '
$ws.Range("E30").Value = '  -3.64%  '
$ws.Range("D31").Value = "'0.117"
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("E32").Value = '  +1.89%  '
$ws.Range("D33").Value = "'4.66"
$ws.Range("E33").Value = '  -3.66%  '
$ws.Range("D34").Value = "'0.0626"
$ws.Range("E34").Value = '  +2.47%  '
$ws.Range("D35").Value = "'4.27"
$ws.Range("E35").Value = '  -1.95%  '
$ws.Range("D36").Value = "'6.21"
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("E38").Value = '  -2.96%  '
$ws.Range("D39").Value = "'2.13"
$ws.Range("E39").Value = '  -6.67%  '
$ws.Range("E40").Value = '  -2.29%  '
$ws.Range("D41").Value = "'0.0973"
$ws.Range("E41").Value = '  -0.58%  '
$ws.Range("E42").Value = '  +4.86%  '
$ws.Range("D43").Value = "'1.18"
$ws.Range("E43").Value = '  -3.04%  '
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("D45").Value = "'16.07"
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("D46").Value = "'1.351.55"
$ws.Range("E46").Value = '  +1.31%  '
$ws.Range("E47").Value = '  -4.73%  '
$ws.Range("D48").Value = "'87.60"
$ws.Range("E48").Value = '  -4.74%  '
$ws.Range("E49").Value = '  -4.40%  '
$ws.Range("D50").Value = "'2.82"
$ws.Range("E50").Value = '  +3.37%  '
$ws.Range("D51").Value = "'45.38"
$ws.Range("E51").Value = '  +4.07%  '
